$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain decimal number (e.g. "304.16").
# Assigning such a string directly would make Excel auto-convert it to a numeric
# value, losing the original text formatting (e.g. trailing zeros). Temporarily mark
# these cells as Text, write the values, then restore the original "General" format.
$textForceRows = 5,6,7,9,10,11,13,14,17,19,22,23,26,27,28,29,30,31,33,36,44,45,47,48,50
foreach ($r in $textForceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '42.999.02'
$ws.Range("E2").Value = '  +1.84%  '
$ws.Range("D3").Value = '2.307.60'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '304.16'
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("D6").Value = '100.44'
$ws.Range("E6").Value = '  +5.64%  '
$ws.Range("D7").Value = '0.504'
$ws.Range("E7").Value = '  +2.25%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '0.515'
$ws.Range("E9").Value = '  +4.38%  '
$ws.Range("D10").Value = '34.97'
$ws.Range("E10").Value = '  +4.83%  '
$ws.Range("D11").Value = '0.0797'
$ws.Range("E11").Value = '  +0.86%  '
$ws.Range("E12").Value = '  +4.19%  '
$ws.Range("D13").Value = '17.99'
$ws.Range("E13").Value = '  +15.95%  '
$ws.Range("D14").Value = '6.89'
$ws.Range("E14").Value = '  +3.42%  '
$ws.Range("D15").Value = '2.684.42'
$ws.Range("E15").Value = '  +2.35%  '
$ws.Range("D16").Value = '2.303.05'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '0.815'
$ws.Range("E17").Value = '  +4.42%  '
$ws.Range("D18").Value = '42.921.39'
$ws.Range("E18").Value = '  +1.87%  '
$ws.Range("D19").Value = '12.48'
$ws.Range("E19").Value = '  +6.78%  '
$ws.Range("E20").Value = '  +2.51%  '
$ws.Range("E21").Value = '  +1.27%  '
$ws.Range("D22").Value = '67.71'
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("D23").Value = '237.18'
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("E24").Value = '  +12.65%  '
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").Value = '24.80'
$ws.Range("E27").Value = '  +3.75%  '
$ws.Range("D28").Value = '2.30'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '167.71'
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").Value = '33.98'
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("D31").Value = '9.17'
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").Value = '5.01'
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("E34").Value = '  +1.82%  '
$ws.Range("E35").Value = '  +3.98%  '
$ws.Range("D36").Value = '17.02'
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("E38").Value = '  +3.49%  '
$ws.Range("E39").Value = '  +1.45%  '
$ws.Range("E40").Value = '  +3.70%  '
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("E42").Value = '  -6.05%  '
$ws.Range("D43").Value = '2.003.18'
$ws.Range("E43").Value = '  +2.17%  '
$ws.Range("D44").Value = '0.0286'
$ws.Range("E44").Value = '  +2.93%  '
$ws.Range("D45").Value = '10.20'
$ws.Range("E45").Value = '  +6.91%  '
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("D47").Value = '2.84'
$ws.Range("E47").Value = '  +1.77%  '
$ws.Range("D48").Value = '55.42'
$ws.Range("E48").Value = '  +6.30%  '
$ws.Range("D49").Value = '2.528.05'
$ws.Range("E49").Value = '  +1.32%  '
$ws.Range("D50").Value = '1.54'
$ws.Range("E50").Value = '  +4.68%  '
$ws.Range("E51").Value = '  +0.93%  '

# Restore the original "General" number format on those cells.
foreach ($r in $textForceRows) {
    $ws.Range("D$r").NumberFormat = "General"
}
